$wb = $excel.ActiveWorkbook

# --- Rename the existing sheet and add the two new ones, in order ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "logit"

$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "gbtree"

$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws2)
$ws3.Name = "rf"

# --- Classification-report / confusion-matrix text blocks ---
$logitReport = "              precision    recall  f1-score   support`n           0       0.84      0.89      0.86      1300`n           1       0.61      0.52      0.56       458`n    accuracy                           0.79      1758`n   macro avg       0.73      0.70      0.71      1758`nweighted avg       0.78      0.79      0.78      1758`n"
$logitMatrix = "[[1152  148]`n [ 222  236]]"

$gbtreeReport = "              precision    recall  f1-score   support`n           0       0.83      0.90      0.86      1300`n           1       0.63      0.48      0.55       458`n    accuracy                           0.79      1758`n   macro avg       0.73      0.69      0.70      1758`nweighted avg       0.78      0.79      0.78      1758`n"
$gbtreeMatrix = "[[1167  133]`n [ 236  222]]"

$rfReport = "              precision    recall  f1-score   support`n           0       0.79      0.96      0.87      1300`n           1       0.73      0.29      0.42       458`n    accuracy                           0.79      1758`n   macro avg       0.76      0.63      0.64      1758`nweighted avg       0.78      0.79      0.75      1758`n"
$rfMatrix = "[[1249   51]`n [ 323  135]]"

# --- Build the bold / thin-bordered / centered-top header style exactly
#     once (on a scratch cell far outside the used range) so the
#     stylesheet only grows by a single cellXfs record, then fan that
#     style out to every sheet via copy/paste-format. ---
$tmpl = $ws1.Range("Z100")
$tmpl.Value = "x"
$tmpl.Font.Bold = $true
$tmpl.HorizontalAlignment = -4108   # xlCenter
$tmpl.VerticalAlignment = -4160     # xlTop
$tmpl.Borders.LineStyle = 1         # xlContinuous
$tmpl.Copy()

function Fill-ModelSheet($ws, $accuracy, $report, $matrix) {
    $ws.Range("B1").Value = "accuracy"
    $ws.Range("C1").Value = "classification report"
    $ws.Range("D1").Value = "confusion materix"
    $ws.Range("A2").Value = 0
    $ws.Range("B2").Value = $accuracy
    $ws.Range("C2").Value = $report
    $ws.Range("D2").Value = $matrix

    $ws.Range("B1:D1").PasteSpecial(-4122)  # xlPasteFormats
    $ws.Range("A2").PasteSpecial(-4122)     # xlPasteFormats

    # Undo the autofit row-height bump the multi-line report/matrix text
    # triggers, so row 2 keeps the sheet's default (no explicit height).
    $ws.Rows.Item(2).AutoFit()
}

Fill-ModelSheet $ws1 0.7895335608646189 $logitReport $logitMatrix
Fill-ModelSheet $ws2 0.7901023890784983 $gbtreeReport $gbtreeMatrix
Fill-ModelSheet $ws3 0.7872582480091013 $rfReport $rfMatrix

$tmpl.Clear()
$ws1.Activate()
